# Vikor / mcdm.xlsx - "Experiments of Case Study Finished"
# Updates the weight row (E5:J5) and the alternative rows (E7:J9) on Sheet1
# with the final experiment results, and moves the active selection to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5: weight vector -------------------------------------------------
$ws.Range("E5").Value = 0.088525943365128895
$ws.Range("F5").Value = -0.0014589542245943799
$ws.Range("G5").Value = 0.049116053092724599
$ws.Range("H5").Value = 0.0031489210752849899
$ws.Range("I5").Value = -0.0016964319232541101
$ws.Range("J5").Value = 0.85605369631901296

# --- Row 7: SPEA-II results (only the VIKOR/Q score column changes) -------
$ws.Range("J7").Value = 0.01363

# --- Row 8: NSGA-II results (only the VIKOR/Q score column changes) -------
$ws.Range("J8").Value = 0

# --- Row 9: Hybrid results (all metric columns change) --------------------
$ws.Range("E9").Value = 96.9
$ws.Range("F9").Value = 0.74052000000000007
$ws.Range("G9").Value = 4606.7493699999995
$ws.Range("H9").Value = 0.20927000000000001
$ws.Range("I9").Value = 133.47863999999998
$ws.Range("J9").Value = 1

# --- Selection moves from J5 to D5 -----------------------------------------
$null = $ws.Activate()
$null = $ws.Range("D5").Select()

Write-Output "edits applied"
